$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6229065656661987
$ws.Range("B1").Value = 2.863307237625122
$ws.Range("C1").Value = 3.109352111816406
$ws.Range("D1").Value = 3.717393398284912
$ws.Range("E1").Value = 1.360589027404785
